# Slide 2: the old "LetterBox" rectangle (Rectangle 1) is removed and its
# letter text is moved onto "Rectangle 2" (which is resized/repositioned to
# the new frame and keeps the border it already had).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Find the two shapes by name so the script doesn't depend on a fixed index.
$letterBox = $null
$rectangle2 = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 1") { $letterBox = $sh }
    if ($sh.Name -eq "Rectangle 2") { $rectangle2 = $sh }
}

# Drop the now-unused LetterBox rectangle entirely.
if ($letterBox -ne $null) {
    $letterBox.Delete()
}

# Move/resize Rectangle 2 onto the new frame (values are EMU / 12700 = points).
$rectangle2.Left = 2126774 / 12700
$rectangle2.Top = 567559 / 12700
$rectangle2.Width = 8120970 / 12700
$rectangle2.Height = 4921338 / 12700

# Give Rectangle 2 the letter text that used to live in the LetterBox shape.
$tr = $rectangle2.TextFrame.TextRange
$tr.Text = "a"
$tr.Font.Size = 413
$tr.Font.Name = "lucida grande"
$tr.Font.Color.RGB = 0x222222
